# Rename the worksheet "Data_Final" to "Data-Final" to avoid loading
# errors when the workbook is read from Jupyter (underscore -> hyphen).
$wb = $excel.ActiveWorkbook

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "Data_Final") {
        $ws = $sheet
        break
    }
}
if ($ws -eq $null) {
    $ws = $wb.Worksheets.Item(1)
}

$ws.Name = "Data-Final"
